# Apply the edit described by the diff: append BIDS dataset metadata rows
# to the "Tabelle2" worksheet, and make that sheet the active / selected one.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Tabelle2")

# New metadata rows appended below the existing Name / BIDSVersion rows.
$ws2.Range("A3").Value = "Name"
$ws2.Range("B3").Value = "Test_template"

$ws2.Range("A4").Value = "DatasetType"
$ws2.Range("B4").Value = "raw"

$ws2.Range("A5").Value = "Licence"
$ws2.Range("B5").Value = "CCO"

$ws2.Range("A6").Value = "Authors"
$ws2.Range("B6").Value = "Karl Koschutnig; Max Mustermann"

$ws2.Range("A7").Value = "Funding"

# Activate Tabelle2 and select the cell shown in the saved selection, mirroring
# the resulting tabSelected / activeTab / selection state in the XML.
$ws2.Activate()
$ws2.Range("C12").Select()

# Match the page setup (paper size / orientation) recorded for the sheet.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
